# Insert two new rows of daily price data for "Ajo" (Chino, Primera)
# at the top of the existing data block (right after the header row),
# pushing the previously-existing rows (old 1161..1255) down to
# (1163..1257), and fill the two new rows with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 1161 (shifts data down).
$ws.Range("A1161:R1162").EntireRow.Insert()

# --- New row 1161 ---
$ws.Range("A1161").Value2 = 10
$ws.Range("B1161").Value2 = "Vega Modelo de Temuco"
$ws.Range("C1161").Value2 = "La Araucanía"
$ws.Range("D1161").Value2 = 45223
$ws.Range("E1161").Value2 = 9
$ws.Range("F1161").Value2 = 100112003
$ws.Range("G1161").Value2 = "Ajo"
$ws.Range("H1161").Value2 = "Chino"
$ws.Range("I1161").Value2 = "Primera"
$ws.Range("J1161").Value2 = 230
$ws.Range("K1161").Value2 = 22000
$ws.Range("L1161").Value2 = 23000
$ws.Range("M1161").Value2 = 22652
$ws.Range("N1161").Value2 = "`$/caja 10 kilos"
$ws.Range("O1161").Value2 = "China"
$ws.Range("P1161").Value2 = 2265
$ws.Range("Q1161").Value2 = 10
$ws.Range("R1161").Value2 = "Hortaliza"

# --- New row 1162 ---
$ws.Range("A1162").Value2 = 10
$ws.Range("B1162").Value2 = "Vega Modelo de Temuco"
$ws.Range("C1162").Value2 = "La Araucanía"
$ws.Range("D1162").Value2 = 45223
$ws.Range("E1162").Value2 = 9
$ws.Range("F1162").Value2 = 100112003
$ws.Range("G1162").Value2 = "Ajo"
$ws.Range("H1162").Value2 = "Chino"
$ws.Range("I1162").Value2 = "Primera"
$ws.Range("J1162").Value2 = 180
$ws.Range("K1162").Value2 = 25000
$ws.Range("L1162").Value2 = 26000
$ws.Range("M1162").Value2 = 25556
$ws.Range("N1162").Value2 = "`$/malla 10 kilos"
$ws.Range("O1162").Value2 = "China"
$ws.Range("P1162").Value2 = 2556
$ws.Range("Q1162").Value2 = 10
$ws.Range("R1162").Value2 = "Hortaliza"
